$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: Professional summary paragraph - neutralize "all Black and
# Asian-American voters" to "50M voters".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "affecting all Black and Asian-American voters, developed geospatial ML",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "affecting 50M voters, developed geospatial ML", 2)

# ---------------------------------------------------------------------------
# Change 2: Siege Analytics bullet - split the run so the new "50M" figure
# is bold / colored like the other highlighted stats in that bullet.
# ---------------------------------------------------------------------------
$afterRng = $d.Content
$afterRng.Find.Execute("race coding errors affecting", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$afterPos = $afterRng.End

$targetRng = $d.Range($afterPos, $d.Content.End)
$targetRng.Find.Execute("all Black and Asian-American", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$targetRng.Text = "50M"
$targetRng.Font.Bold = $true
$targetRng.Font.Color = 5258796

# ---------------------------------------------------------------------------
# Change 3: Reorder PROFESSIONAL EXPERIENCE entries - move the
# "Analytics Supervisor - GSD&M" block (heading + 4 bullets) so it comes
# after the "Data Products Manager - Helm/Murmuration" block instead of
# before it.
# ---------------------------------------------------------------------------
for ($i = 0; $i -lt 5; $i++) {
    $firstOfBlock = $d.Paragraphs(15)
    $moveRange = $d.Range($firstOfBlock.Range.Start, $firstOfBlock.Range.End)
    $moveRange.Cut()

    $anchorRng = $d.Content
    $anchorRng.Find.Execute("Senior Analyst - Myers Research", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $insertPos = $anchorRng.Start

    $pasteRange = $d.Range($insertPos, $insertPos)
    $pasteRange.Paste()
}

# ---------------------------------------------------------------------------
# Change 4: Key Projects impact statement - neutralize language.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Impact: Corrected demographic data affecting all Black and Asian-American voters, improved",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Impact: Corrected demographic data affecting 50M voters nationwide, improved", 2)

Write-Output "edit complete"
